$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Word's InlineShape object has no settable .Name — the only way to
# rewrite the underlying docPr/cNvPr "name" (the picture's internal
# filename label) through the object model is to flip the picture to a
# floating Shape (which does expose .Name), rename it, then flip it
# back to an inline picture.
function Rename-InlinePicture($headerFooter, $newName) {
    $ishp = $headerFooter.Range.InlineShapes(1)
    $shp = $ishp.ConvertToShape()
    $shp.Name = $newName
    $shp.ConvertToInlineShape() | Out-Null
}

# Primary (default) footer -> footer2.xml, Pearson logo id="4"
Rename-InlinePicture $sec.Footers(1) "image2.png"

# First-page footer -> footer1.xml, Pearson logo id="2"
Rename-InlinePicture $sec.Footers(2) "image2.png"

# First-page header -> header1.xml, BTEC logo id="1"
Rename-InlinePicture $sec.Headers(2) "image1.jpg"

# Primary (default) header -> header2.xml, BTEC logo id="3"
Rename-InlinePicture $sec.Headers(1) "image1.jpg"
